# cap nhat thong ke
# Update the "Diem tong ket" (summary score) statistics sheet:
#  - Group N7 (row 8) now has a final score of 400/400 (was blank/0) and
#    total 8 (was 0).
#  - Group N17 (row 18) final score corrected from 206/400 to 400/400,
#    and total corrected from 5 to 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "400/400"
$ws.Range("D8").Value = 8

$ws.Range("B18").Value = "400/400"
$ws.Range("D18").Value = 8

# Update the last active selection shown in the sheet view.
[void]$ws.Range("I11").Select()
